$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column J: "Daily First Doses" = day-over-day change in I (Total First Doses) ---
$ws.Range("J1").Value = "Daily First Doses"

# --- New data row 112 (day 110 / 2021-04-09) ---
$ws.Range("A112").Value = 110
$ws.Range("B112").Value = 44304
$ws.Range("C112").Value = 86565
$ws.Range("D112").Value = 3837881
$ws.Range("E112").Formula = "=(D112-F112)"
$ws.Range("F112").Value = 690620
$ws.Range("G112").Value = 345310
$ws.Range("H112").Formula = "=AVERAGE(C106:C112)"
$ws.Range("I112").Formula = "=(D112-G112)"

# --- Fill the new "Daily First Doses" formula down J10:J112 ---
for ($r = 10; $r -le 112; $r++) {
    $prev = $r - 1
    $ws.Range("J$r").Formula = "=(I$r-I$prev)"
}

# --- Column width for the new column ---
$ws.Columns.Item(10).ColumnWidth = 14.6640625

# --- View state (scrolled down, new active cell selected) ---
$ws.Application.ActiveWindow.ScrollRow = 90
$ws.Range("F114").Select()
